# Add a new "Tone detection threshold" column to the summary sheet.
#
# The sheet has a 3-row header (rows 1-3) followed by per-subject data rows
# (4-23). Columns K:P hold the "Speech IM (experiment 2)" block (Noise /
# Speech / difference(dB) / left / right / average). We insert one new
# column before K to hold a "Tone detection threshold" value per subject,
# which pushes the existing K:P block to L:Q (and the dimension/mergeCells
# follow automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at K - this shifts the existing K:P data/styles to
#    L:Q, extends merged cells (K1:M2 -> L1:N2, N1:P2 -> O1:Q2) and the used
#    range/dimension to column Q automatically.
$ws.Columns("K").Insert()

# 2) The inserted column copies formatting from its left neighbour (J), so
#    re-stamp K1:K3 with the "header band" look already used elsewhere on
#    the sheet (light red/rose fill, thin border, centered) - e.g. J2 /
#    B3:J3 already carry that exact style. Copy format only (not values).
$ws.Range("J2").Copy()
$ws.Range("K1:K3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# 3) New column header/label.
$ws.Range("K3").Value = "Tone detection threshold"

# 4) Per-subject tone detection threshold values (row 4 = TLZ ... row 23 = TNM).
$ws.Range("K4").Value = 5
$ws.Range("K5").Value = 5
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 5
$ws.Range("K8").Value = 5
$ws.Range("K9").Value = 5
$ws.Range("K10").Value = 5
$ws.Range("K11").Value = 7.5
$ws.Range("K12").Value = 7.5
$ws.Range("K13").Value = 10
$ws.Range("K14").Value = 5
$ws.Range("K15").Value = 10
$ws.Range("K16").Value = 5
$ws.Range("K17").Value = 5
$ws.Range("K18").Value = 10
$ws.Range("K19").Value = 5
$ws.Range("K20").Value = 12.5
$ws.Range("K21").Value = 5
$ws.Range("K22").Value = 5
$ws.Range("K23").Value = 20

# 5) Column widths: K becomes the wide label column, J (the column right
#    before it) also widens slightly.
$ws.Columns("J").ColumnWidth = 10.67
$ws.Columns("K").ColumnWidth = 26.83
